# Updates the cryptocurrency price/volume table (columns B-E, rows 2-51)
# to match the refreshed data from the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to hold a literal text value (not an auto-converted
    # number/date), then restore the original (default) cell style so no
    # formatting side effects are introduced.
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '42.237.89'
$ws.Range("E2").Value = '  +0.45%  '

$ws.Range("D3").Value = '2.282.86'
$ws.Range("E3").Value = '  -0.90%  '

$ws.Range("E4").Value = '  +0.02%  '

Set-TextValue "D5" '323.39'
$ws.Range("E5").Value = '  +2.35%  '

Set-TextValue "D6" '101.78'
$ws.Range("E6").Value = '  -2.60%  '

$ws.Range("E7").Value = '  -0.29%  '

$ws.Range("E8").Value = '  +0.08%  '

Set-TextValue "D9" '0.606'
$ws.Range("E9").Value = '  -0.32%  '

Set-TextValue "D10" '39.70'
$ws.Range("E10").Value = '  -0.05%  '

Set-TextValue "D11" '0.0901'
$ws.Range("E11").Value = '  -0.62%  '

Set-TextValue "D12" '8.32'
$ws.Range("E12").Value = '  -1.62%  '

$ws.Range("E13").Value = '  -0.96%  '

Set-TextValue "D14" '0.964'
$ws.Range("E14").Value = '  -0.99%  '

Set-TextValue "D15" '15.09'
$ws.Range("E15").Value = '  -2.37%  '

$ws.Range("D16").Value = '2.629.45'
$ws.Range("E16").Value = '  -0.82%  '

$ws.Range("D17").Value = '2.281.55'
$ws.Range("E17").Value = '  -0.94%  '

$ws.Range("D18").Value = '42.227.59'
$ws.Range("E18").Value = '  +0.50%  '

Set-TextValue "D19" '7.37'
$ws.Range("E19").Value = '  -4.49%  '

$ws.Range("E20").Value = '  -0.21%  '

Set-TextValue "D21" '12.77'
$ws.Range("E21").Value = '  +27.62%  '

Set-TextValue "D22" '3.66'

Set-TextValue "D23" '73.07'
$ws.Range("E23").Value = '  -1.15%  '

Set-TextValue "D24" '268.54'
$ws.Range("E24").Value = '  -6.58%  '

$ws.Range("E25").Value = '  -3.50%  '

$ws.Range("E26").Value = '  -0.19%  '

Set-TextValue "D27" '10.85'
$ws.Range("E27").Value = '  -0.79%  '

Set-TextValue "D28" '2.29'
$ws.Range("E28").Value = '  +2.89%  '

Set-TextValue "D29" '22.42'
$ws.Range("E29").Value = '  -4.87%  '

Set-TextValue "D30" '38.05'
$ws.Range("E30").Value = '  +6.98%  '

Set-TextValue "D31" '164.12'
$ws.Range("E31").Value = '  -0.88%  '

Set-TextValue "D32" '6.09'
$ws.Range("E32").Value = '  +3.25%  '

Set-TextValue "D33" '0.0877'
$ws.Range("E33").Value = '  -0.68%  '

$ws.Range("E34").Value = '  +1.04%  '

$ws.Range("E35").Value = '  -13.05%  '

$ws.Range("E36").Value = '  -3.39%  '

$ws.Range("E37").Value = '  -1.22%  '

Set-TextValue "D38" '0.0353'
$ws.Range("E38").Value = '  +0.20%  '

Set-TextValue "D39" '3.70'
$ws.Range("E39").Value = '  +2.08%  '

Set-TextValue "D40" '2.75'
$ws.Range("E40").Value = '  -6.59%  '

$ws.Range("E41").Value = '  +1.82%  '

Set-TextValue "D42" '68.85'
$ws.Range("E42").Value = '  -2.63%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue "D43" '0.225'
$ws.Range("E43").Value = '  -0.65%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D44" '1.00'
$ws.Range("E44").Value = '  -0.06%  '

Set-TextValue "D45" '90.80'
$ws.Range("E45").Value = '  -12.29%  '

Set-TextValue "D46" '12.20'
$ws.Range("E46").Value = '  +1.17%  '

Set-TextValue "D47" '112.94'
$ws.Range("E47").Value = '  -3.55%  '

Set-TextValue "D48" '79.98'
$ws.Range("E48").Value = '  +3.00%  '

Set-TextValue "D49" '8.92'
$ws.Range("E49").Value = '  -2.00%  '

$ws.Range("E50").Value = '  -2.14%  '

$ws.Range("D51").Value = '1.590.38'
$ws.Range("E51").Value = '  +1.86%  '
